$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values would otherwise be auto-coerced to Number by Excel
# (losing a significant trailing zero), so force Text format first.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.494.86"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.678.50"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "219.88"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = "0.5316"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.2702"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").Value = "0.06412"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "21.86"
$ws.Range("E10").Value = "  +5.37%  "
$ws.Range("D11").Value = "0.07800"
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.693.05"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.515"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("D14").Value = "0.5599"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "0.0₅8359"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "65.73"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "26.527.31"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "4.802"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").Value = "193.49"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "10.32"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "6.341"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "0.1277"
$ws.Range("E24").Value = "  +6.01%  "
$ws.Range("D25").Value = "139.50"
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("D26").Value = "7.423"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "1.443"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("D29").Value = "0.06288"
$ws.Range("E29").Value = "  +6.98%  "
$ws.Range("D30").Value = "1.288"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("D31").Value = "3.617"
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("D32").Value = "3.461"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").Value = "1.699"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").Value = "1.015"
$ws.Range("E34").Value = "  +3.14%  "
$ws.Range("D35").Value = "0.6181"
$ws.Range("E35").Value = "  +9.16%  "
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").Value = "2.788"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").Value = "6.159"
$ws.Range("E38").Value = "  +7.82%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "1.095.60"
$ws.Range("E40").Value = "  +6.55%  "
$ws.Range("D41").Value = "0.8619"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "0.9997"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "100.70"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").Value = "1.821.77"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").Value = "58.80"
$ws.Range("E45").Value = "  +5.31%  "
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "8.187"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.484"
$ws.Range("E49").Value = "  +7.31%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05195"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "6.040"
$ws.Range("E51").Value = "  +2.37%  "
